$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Column C ("Förändrad") holds a date serial that was bumped by one day
# (45202 -> 45203) for every data row (rows 2 through 319).
$ws.Range("C2:C319").Value = 45203
